# Actualización 10 de Mayo
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja "Estadisticos 1P"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D3").Value = 16
$ws1.Range("F3").Value = 23
$ws1.Range("G3").Value = 58.97

$ws1.Range("D4").Value = 15
$ws1.Range("F4").Value = 24
$ws1.Range("G4").Value = 61.54
$ws1.Range("H4").Value = 8.7

# ---------------------------------------------------------------------------
# Hoja "Estadisticos 2P"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 17
$ws2.Range("E2").Value = 12
$ws2.Range("F2").Value = 23
$ws2.Range("G2").Value = 57.5

$ws2.Range("D3").Value = 26
$ws2.Range("E3").Value = 10
$ws2.Range("F3").Value = 13
$ws2.Range("G3").Value = 33.33
$ws2.Range("H3").Value = 9.8

$ws2.Range("D4").Value = 20
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Value = 19
$ws2.Range("G4").Value = 48.72
$ws2.Range("H4").Value = 8.7

# ---------------------------------------------------------------------------
# Hoja "Estadisticos Final"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D3").Value = 16
$ws3.Range("F3").Value = 23
$ws3.Range("G3").Value = 58.97

$ws3.Range("D4").Value = 15
$ws3.Range("F4").Value = 24
$ws3.Range("G4").Value = 61.54
$ws3.Range("H4").Value = 9

# ---------------------------------------------------------------------------
# Hoja "Rescatables" - se agregan 11 alumnos rescatables
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$nc        = @(20330051920116, 20330051920121, 20330051920389, 20330051920144, 19330051920366, 19330051920375, 19330051920377, 19330051920382, 19330051920362, 19330051920368, 19330051920378)
$paterno   = @("CARRERA", "CUATRA", "PACHECO", "DE LA ROSA", "CRISTOBAL", "HERNANDEZ", "HERNANDEZ", "MAZAHUA", "CALIHUA", "DE LA CRUZ", "HERNANDEZ")
$materno   = @("ROMANOS", "ZOPIYACTLE", "MAZAHUA", "CASTRO", "BRUNO", "ANTONIO", "FLORES", "IXMATLAHUA", "CALIHUA", "DE LA CRUZ", "HERNANDEZ")
$nombres   = @("AMARANTA DENISSE", "MARIA", "TAILY", "ALONDRA", "DANIELA", "MARIA GUADALUPE", "PERLA", "LUCERO", "JOEL", "OFELIA", "DARIANA MONSERRAT")
$nomLargo  = @(
    "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",
    "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",
    "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",
    "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",
    "EVALÚA EL DESEMPEÑO DE LA ORGANIZACIÓN UTILIZANDO HERRAMIENTAS DE CALIDAD",
    "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO",
    "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO",
    "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO",
    "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO",
    "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO",
    "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO"
)
$grupo     = @("2ARHV", "2ARHV", "2ARHV", "2ARHV", "4ARHV", "4ARHV", "4ARHV", "4ARHV", "4ARHV", "4ARHV", "4ARHV")
$reprob    = @(2, 2, 2, 2, 2, 2, 2, 2, 1, 1, 1)

# Se llena columna por columna (B, C, D, E, F, G, A) replicando el orden en
# que las cadenas nuevas fueron agregadas a sharedStrings.xml.
for ($i = 0; $i -lt $paterno.Length; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt $materno.Length; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt $nombres.Length; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt $nomLargo.Length; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 5).Value = $nomLargo[$i]
}
for ($i = 0; $i -lt $grupo.Length; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 6).Value = $grupo[$i]
}
for ($i = 0; $i -lt $reprob.Length; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 7).Value = $reprob[$i]
}
for ($i = 0; $i -lt $nc.Length; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 1).Value = $nc[$i]
}
